# Add new sub2sector "18_01_02_gas_power_ccs".
#
# The reference table lists one row per sub2sector in columns B/C (plus a
# constant "(new)" marker in column D). The three former CCGT breakdown rows
# (18_01_02_01_gasturbine, 18_01_02_02_combinedcycle, 18_01_02_03_ccs) are
# consolidated: the first of the three rows is turned into the new
# "18_01_02_gas_power_ccs" row, and the other two are removed, so every
# following row shifts up by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 247 held "18_01_02_01_gasturbine" in both B and C; repurpose it for the
# newly introduced sub2sector.
$ws.Cells.Item(247, 2).Value2 = "18_01_02_gas_power_ccs"
$ws.Cells.Item(247, 3).Value2 = "18_01_02_gas_power_ccs"

# Rows 248 ("18_01_02_02_combinedcycle") and 249 ("18_01_02_03_ccs") are no
# longer needed; delete them outright so everything below shifts up by two
# rows (old row 250 becomes the new row 248, etc.).
$ws.Range("248:249").Delete() | Out-Null

# Match the author's final selection/view state.
$ws.Range("D247").Select()
